# Weekly update: a new daily record is inserted at row 60 (pushing the
# existing historical rows down by one), matching the commit
# "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 60; this pushes the former rows 60-127
# down to rows 61-128 (dimension grows from R127 to R128 automatically).
$ws.Rows("60:60").Insert()

# Populate the new row 60 with the new daily price record.
$ws.Range("A60").Value = 4
$ws.Range("B60").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C60").Value = "Los Lagos"
$ws.Range("D60").Value = 44494
$ws.Range("E60").Value = 10
$ws.Range("F60").Value = 100112039
$ws.Range("G60").Value = "Ciboulette"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 80
$ws.Range("K60").Value = 2500
$ws.Range("L60").Value = 2500
$ws.Range("M60").Value = 2500
$ws.Range("N60").Value = "$/docena de atados"
$ws.Range("O60").Value = "Región Metropolitana"
$ws.Range("P60").Value = 833
$ws.Range("Q60").Value = 3
$ws.Range("R60").Value = "Hortaliza"
